$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B2").Value = "serie1"
$ws.Range("C2").Value = "serie2"
$ws.Range("D2").Value = "serie3"
$ws.Range("E2").Value = "serie4"
$ws.Range("F2").Value = "serie5"
$ws.Range("G2").Value = "serie6"

$ws.Range("B3").Value = 0.10790991783142
$ws.Range("C3").Value = 0.71073627471923795
$ws.Range("D3").Value = 4.5886371135711599
$ws.Range("E3").Value = 31.6078133583068
$ws.Range("F3").Value = 150.43901991844101
$ws.Range("G3").Value = 1137.9384648799801
$ws.Range("B4").Value = 0.11775851249694801
$ws.Range("C4").Value = 0.81680727005004805
$ws.Range("D4").Value = 5.0179984569549498
$ws.Range("E4").Value = 30.526545763015701
$ws.Range("F4").Value = 147.95070028305
$ws.Range("G4").Value = 1086.4372854232699
$ws.Range("B5").Value = 0.115154981613159
$ws.Range("C5").Value = 0.74224162101745605
$ws.Range("D5").Value = 4.8287572860717702
$ws.Range("E5").Value = 22.909975528716998
$ws.Range("F5").Value = 150.113550901412
$ws.Range("G5").Value = 1089.29955005645
$ws.Range("B6").Value = 0.083296060562133706
$ws.Range("C6").Value = 0.71216511726379395
$ws.Range("D6").Value = 5.0334162712097097
$ws.Range("E6").Value = 23.281772375106801
$ws.Range("F6").Value = 148.559745311737
$ws.Range("G6").Value = 1053.6332614421799
$ws.Range("B7").Value = 0.106697559356689
$ws.Range("C7").Value = 0.68281650543212802
$ws.Range("D7").Value = 5.0962104797363201
$ws.Range("E7").Value = 21.350227355956999
$ws.Range("F7").Value = 148.79026055335899
$ws.Range("G7").Value = 1029.3666572570801
$ws.Range("B8").Value = 0.10542106628417899
$ws.Range("C8").Value = 0.71890616416931097
$ws.Range("D8").Value = 4.7525465488433802
$ws.Range("E8").Value = 21.0414991378784
$ws.Range("F8").Value = 147.23434805869999
$ws.Range("G8").Value = 696.00749921798695
$ws.Range("B9").Value = 0.077968597412109306
$ws.Range("C9").Value = 0.67446494102478005
$ws.Range("D9").Value = 4.8188965320587096
$ws.Range("E9").Value = 21.358883619308401
$ws.Range("F9").Value = 148.84608340263301
$ws.Range("G9").Value = 900.32061719894398
$ws.Range("B10").Value = 0.101251125335693
$ws.Range("C10").Value = 0.53674960136413497
$ws.Range("D10").Value = 5.3744482994079501
$ws.Range("E10").Value = 20.965662240981999
$ws.Range("F10").Value = 147.84991526603699
$ws.Range("G10").Value = 828.62716650962795
$ws.Range("B11").Value = 0.119367837905883
$ws.Range("C11").Value = 0.72379612922668402
$ws.Range("D11").Value = 6.6076405048370299
$ws.Range("E11").Value = 21.2100734710693
$ws.Range("F11").Value = 147.50043892860401
$ws.Range("G11").Value = 708.21263694763104
$ws.Range("B12").Value = 0.099427223205566406
$ws.Range("C12").Value = 0.72187829017639105
$ws.Range("D12").Value = 5.3477084636688197
$ws.Range("E12").Value = 21.311897754669101
$ws.Range("F12").Value = 163.27137827873199
$ws.Range("G12").Value = 805.32461428642205

$ws.Range("A14").Value = "moy"
$ws.Range("B14").Formula = "=AVERAGE(B3:B12)"
$ws.Range("C14:G14").Formula = "=AVERAGE(C3:C12)"

$ws.Range("G11").Select()
